$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Toggle several "missing" values in column E (imputation corrections) ---
# E5: -5 -> missing
$ws.Range("E5").ClearContents()
# E8: missing -> -6.6
$ws.Range("E8").Value = -6.6
# E12: -5.3 -> missing
$ws.Range("E12").ClearContents()
# E14: missing -> -5.4
$ws.Range("E14").Value = -5.4
# E18: -8.5 -> missing
$ws.Range("E18").ClearContents()

# --- Remove two whole data rows: "RM 232" (row 26) and "SC 92" (row 28) ---
# Delete the higher-numbered row first so the lower row index stays valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# --- After the deletions, rows 27-35 have shifted up to 26-33. ---
# Row 26 is now "SC 5"; its C column value is no longer missing.
$ws.Range("C26").Value = 10.8

# Row 27 is now "SC 101"; its C column value became missing.
$ws.Range("C27").ClearContents()

# Row 33 is now "SC 232"; its D column value is no longer missing.
$ws.Range("D33").Value = -14.1
